# Rename the worksheet tab/sheet name from "Sheet1" to "Nathan Ellis"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Nathan Ellis"

# Insert a new column before column A to make room for "matchNo"
$ws.Columns.Item(1).Insert()

# Insert a new row before row 2 to make room for the "37th" match entry
$ws.Rows.Item(2).Insert()

# --- Row 1: headers ---
$ws.Cells.Item(1, 1).Value = "matchNo"
$ws.Cells.Item(1, 2).Value = "teamName"
$ws.Cells.Item(1, 3).Value = "batterName"
$ws.Cells.Item(1, 4).Value = "states"
$ws.Cells.Item(1, 5).Value = "runs"
$ws.Cells.Item(1, 6).Value = "balls"
$ws.Cells.Item(1, 7).Value = "fours"
$ws.Cells.Item(1, 8).Value = "sixes"
$ws.Cells.Item(1, 9).Value = "sr"
$ws.Cells.Item(1, 10).Value = "opponentTeamName"
$ws.Cells.Item(1, 11).Value = "venue"
$ws.Cells.Item(1, 12).Value = "date"
$ws.Cells.Item(1, 13).Value = "result"

# Mark the numeric-looking columns as Text before writing so values like
# "12" / "0" / "100.00" persist as strings (matching the source data's
# t="str" cells) instead of being auto-coerced into numbers.
$ws.Range("E2:I3").NumberFormat = "@"

# --- Row 2: new match ("37th") ---
$ws.Cells.Item(2, 1).Value = "37th"
$ws.Cells.Item(2, 2).Value = "Punjab Kings"
$ws.Cells.Item(2, 3).Value = "Nathan Ellis"
$ws.Cells.Item(2, 4).Value = "c Pandey b Kumar"
$ws.Cells.Item(2, 5).Value = "12"
$ws.Cells.Item(2, 6).Value = "12"
$ws.Cells.Item(2, 7).Value = "0"
$ws.Cells.Item(2, 8).Value = "1"
$ws.Cells.Item(2, 9).Value = "100.00"
$ws.Cells.Item(2, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(2, 11).Value = "Sharjah"
$ws.Cells.Item(2, 12).Value = "September 25"
$ws.Cells.Item(2, 13).Value = "Punjab Kings won by 5 runs"

# --- Row 3: existing match ("42nd"), shifted down with new matchNo column ---
$ws.Cells.Item(3, 1).Value = "42nd"
$ws.Cells.Item(3, 2).Value = "Punjab Kings"
$ws.Cells.Item(3, 3).Value = "Nathan Ellis"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "6"
$ws.Cells.Item(3, 6).Value = "4"
$ws.Cells.Item(3, 7).Value = "0"
$ws.Cells.Item(3, 8).Value = "0"
$ws.Cells.Item(3, 9).Value = "150.00"
$ws.Cells.Item(3, 10).Value = "Mumbai Indians"
$ws.Cells.Item(3, 11).Value = "Abu Dhabi"
$ws.Cells.Item(3, 12).Value = "September 28"
$ws.Cells.Item(3, 13).Value = "Mumbai won by 6 wickets (with 6 balls remaining)"

# Restore the default "Normal" style on those cells so the saved XML has no
# stray style index (keeps the text type but drops the explicit "@" format).
$ws.Range("E2:I3").Style = "Normal"
